# "Test result grouping in real rendering"
#
# Sheet1: the three detail rows that belong to the first customer group
# (rows 3-5) and the two detail rows of the second group (rows 8-9) get
# their label column (B) and value column (E) merged into a single
# vertically-centered, bold cell - the classic "group header" look for a
# rendered report. The merged value cells pick up a new bold + vertically
# centered style.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# New style used by the merged "group value" cells: bold font, vertically
# centered text (numeric format / everything else stays default).
$valueCells1 = $ws1.Range("E3:E5")
$valueCells1.Font.Bold = $true
$valueCells1.VerticalAlignment = -4108   # xlVAlignCenter

$valueCells2 = $ws1.Range("E8:E9")
$valueCells2.Font.Bold = $true
$valueCells2.VerticalAlignment = -4108   # xlVAlignCenter

# Group the label column and the value column for each customer block.
$ws1.Range("B3:B5").Merge()
$ws1.Range("E3:E5").Merge()

$ws1.Range("B8:B9").Merge()
$ws1.Range("E8:E9").Merge()
